$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text so values like "40.830.64" are not
# misinterpreted as numbers, then restore the original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.830.64"
$ws.Range("D3").Value = "2.392.01"
$ws.Range("D5").Value = "313.24"
$ws.Range("D6").Value = "88.05"
$ws.Range("D7").Value = "0.532"
$ws.Range("D9").Value = "0.493"
$ws.Range("D10").Value = "0.0827"
$ws.Range("D11").Value = "31.13"
$ws.Range("D13").Value = "2.762.66"
$ws.Range("D14").Value = "6.66"
$ws.Range("D15").Value = "15.12"
$ws.Range("D16").Value = "2.374.53"
$ws.Range("D17").Value = "0.762"
$ws.Range("D18").Value = "40.741.44"
$ws.Range("D19").Value = "0.0₃0914"
$ws.Range("D20").Value = "6.17"
$ws.Range("D21").Value = "69.76"
$ws.Range("D22").Value = "10.79"
$ws.Range("D23").Value = "238.21"
$ws.Range("D24").Value = "2.65"
$ws.Range("D27").Value = "23.78"
$ws.Range("D29").Value = "9.42"
$ws.Range("D30").Value = "34.09"
$ws.Range("D31").Value = "156.65"
$ws.Range("D33").Value = "5.25"
$ws.Range("D34").Value = "0.0734"
$ws.Range("D35").Value = "2.43"
$ws.Range("D36").Value = "0.114"
$ws.Range("D37").Value = "2.82"
$ws.Range("D38").Value = "15.93"
$ws.Range("D39").Value = "1.75"
$ws.Range("D40").Value = "0.0991"
$ws.Range("D41").Value = "3.84"
$ws.Range("D43").Value = "1.969.60"
$ws.Range("D45").Value = "17.92"
$ws.Range("D46").Value = "2.82"
$ws.Range("D47").Value = "9.31"
$ws.Range("D48").Value = "2.628.35"
$ws.Range("D49").Value = "73.41"
$ws.Range("D50").Value = "93.78"
$ws.Range("D51").Value = "50.56"

$ws.Range("D2:D51").Style = "Normal"

$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("E3").Value = "  -3.39%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("E7").Value = "  -3.79%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -4.80%  "
$ws.Range("E10").Value = "  -4.14%  "
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("E16").Value = "  -5.02%  "
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -3.95%  "
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -4.67%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -5.53%  "
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  -4.13%  "
$ws.Range("E30").Value = "  -5.39%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("E35").Value = "  -6.02%  "
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("E38").Value = "  -7.95%  "
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("E42").Value = "  -7.54%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("E45").Value = "  -7.69%  "
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("E51").Value = "  -3.42%  "
